$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.522.53'
$ws.Range("E2").Value = '  -0.50%  '

$ws.Range("D3").Value = '2.070.90'
$ws.Range("E3").Value = '  -0.16%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.65'
$ws.Range("E5").Value = '  -0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.630'
$ws.Range("E6").Value = '  +1.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.54'
$ws.Range("E8").Value = '  -1.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.390'
$ws.Range("E9").Value = '  -1.04%  '

$ws.Range("E10").Value = '  -0.29%  '

$ws.Range("E11").Value = '  +1.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.96'
$ws.Range("E12").Value = '  +1.39%  '

$ws.Range("D13").Value = '2.375.57'
$ws.Range("E13").Value = '  -0.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.96'
$ws.Range("E14").Value = '  -0.02%  '

$ws.Range("E15").Value = '  -1.60%  '

$ws.Range("E16").Value = '  -0.66%  '

$ws.Range("D17").Value = '2.057.91'
$ws.Range("E17").Value = '  -1.39%  '

$ws.Range("D18").Value = '37.476.15'
$ws.Range("E18").Value = '  -0.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.56'
$ws.Range("E19").Value = '  -1.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.96'
$ws.Range("E20").Value = '  -2.64%  '

$ws.Range("D21").Value = '0.0₃0829'
$ws.Range("E21").Value = '  -0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.21'
$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("E24").Value = '  +0.59%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  -1.33%  '

$ws.Range("E26").Value = '  +6.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.56'
$ws.Range("E27").Value = '  -1.17%  '

$ws.Range("E28").Value = '  -2.88%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.47'
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("E30").Value = '  -0.86%  '

$ws.Range("E31").Value = '  +1.33%  '

$ws.Range("E32").Value = '  -1.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0633'
$ws.Range("E33").Value = '  +0.16%  '

$ws.Range("E34").Value = '  -0.47%  '

$ws.Range("E35").Value = '  -0.17%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.33'
$ws.Range("E37").Value = '  -1.90%  '

$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("E39").Value = '  -1.10%  '

$ws.Range("E40").Value = '  +7.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.69'
$ws.Range("E41").Value = '  +0.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.21'
$ws.Range("E42").Value = '  +4.79%  '

$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0958'
$ws.Range("E43").Value = '  -1.72%  '

$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.91'
$ws.Range("E44").Value = '  +1.24%  '

$ws.Range("D45").Value = '1.480.84'
$ws.Range("E45").Value = '  +3.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.71'
$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("E47").Value = '  -1.53%  '

$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.99'
$ws.Range("E48").Value = '  -4.96%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.26'
$ws.Range("E49").Value = '  -1.68%  '

$ws.Range("E50").Value = '  -2.05%  '

$ws.Range("D51").Value = '2.258.91'
$ws.Range("E51").Value = '  -0.18%  '
